$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, border, centered) from an existing
# header cell onto the three new header cells, then set their text.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every player row.
$lastRow = 55
$ws.Range("AD2:AD$lastRow").Value = 79
$ws.Range("AE2:AE$lastRow").Value = 82
$ws.Range("AF2:AF$lastRow").Value = 0
